$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 09:32:47"
$ws.Cells.Item(3,1).Value = "Total filas: 107"
$ws.Cells.Item(6,1).Value = "04:03:21"
$ws.Cells.Item(6,2).Value = "04:46"
$ws.Cells.Item(6,3).Value = "215A_EL PATO"
$ws.Cells.Item(6,4).Value = 43
$ws.Cells.Item(6,5).Value = "LP1912"
$ws.Cells.Item(7,1).Value = "04:48:48"
$ws.Cells.Item(7,2).Value = "04:53"
$ws.Cells.Item(7,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(7,4).Value = 5
$ws.Cells.Item(7,5).Value = "LP1912"
$ws.Cells.Item(8,1).Value = "04:48:48"
$ws.Cells.Item(8,2).Value = "05:16"
$ws.Cells.Item(8,3).Value = "17_ROMERO"
$ws.Cells.Item(8,4).Value = 28
$ws.Cells.Item(8,5).Value = "LP1912"
$ws.Cells.Item(9,1).Value = "04:48:48"
$ws.Cells.Item(9,2).Value = "05:22"
$ws.Cells.Item(9,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(9,4).Value = 34
$ws.Cells.Item(9,5).Value = "LP1912"
$ws.Cells.Item(10,1).Value = "05:20:44"
$ws.Cells.Item(10,2).Value = "05:23"
$ws.Cells.Item(10,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(10,4).Value = 3
$ws.Cells.Item(10,5).Value = "LP1912"
$ws.Cells.Item(11,1).Value = "05:20:44"
$ws.Cells.Item(11,2).Value = "05:34"
$ws.Cells.Item(11,3).Value = "215B_EL PATO"
$ws.Cells.Item(11,4).Value = 14
$ws.Cells.Item(11,5).Value = "LP1912"
$ws.Cells.Item(12,1).Value = "05:20:44"
$ws.Cells.Item(12,2).Value = "05:46"
$ws.Cells.Item(12,3).Value = "15_ABASTO"
$ws.Cells.Item(12,4).Value = 26
$ws.Cells.Item(12,5).Value = "LP1912"
$ws.Cells.Item(13,1).Value = "05:20:44"
$ws.Cells.Item(13,2).Value = "05:54"
$ws.Cells.Item(13,3).Value = "10_OLMOS"
$ws.Cells.Item(13,4).Value = 34
$ws.Cells.Item(13,5).Value = "LP1912"
$ws.Cells.Item(14,1).Value = "05:20:44"
$ws.Cells.Item(14,2).Value = "06:04"
$ws.Cells.Item(14,3).Value = "16_SANTA ANA"
$ws.Cells.Item(14,4).Value = 44
$ws.Cells.Item(14,5).Value = "LP1912"
$ws.Cells.Item(15,1).Value = "06:03:38"
$ws.Cells.Item(15,2).Value = "06:11"
$ws.Cells.Item(15,3).Value = "215A_EL PATO"
$ws.Cells.Item(15,4).Value = 8
$ws.Cells.Item(15,5).Value = "LP1912"
$ws.Cells.Item(16,1).Value = "06:03:38"
$ws.Cells.Item(16,2).Value = "06:14"
$ws.Cells.Item(16,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(16,4).Value = 11
$ws.Cells.Item(16,5).Value = "LP1912"
$ws.Cells.Item(17,1).Value = "06:03:38"
$ws.Cells.Item(17,2).Value = "06:21"
$ws.Cells.Item(17,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(17,4).Value = 18
$ws.Cells.Item(17,5).Value = "LP1912"
$ws.Cells.Item(18,1).Value = "06:03:38"
$ws.Cells.Item(18,2).Value = "06:27"
$ws.Cells.Item(18,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(18,4).Value = 24
$ws.Cells.Item(18,5).Value = "LP1912"
$ws.Cells.Item(19,1).Value = "06:03:38"
$ws.Cells.Item(19,2).Value = "06:29"
$ws.Cells.Item(19,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(19,4).Value = 26
$ws.Cells.Item(19,5).Value = "LP1912"
$ws.Cells.Item(20,1).Value = "06:03:38"
$ws.Cells.Item(20,2).Value = "06:31"
$ws.Cells.Item(20,3).Value = "16_SANTA ANA"
$ws.Cells.Item(20,4).Value = 28
$ws.Cells.Item(20,5).Value = "LP1912"
$ws.Cells.Item(21,1).Value = "06:03:38"
$ws.Cells.Item(21,2).Value = "06:44"
$ws.Cells.Item(21,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(21,4).Value = 41
$ws.Cells.Item(21,5).Value = "LP1912"
$ws.Cells.Item(22,1).Value = "06:03:38"
$ws.Cells.Item(22,2).Value = "06:46"
$ws.Cells.Item(22,3).Value = "215C_EL PATO"
$ws.Cells.Item(22,4).Value = 43
$ws.Cells.Item(22,5).Value = "LP1912"
$ws.Cells.Item(23,1).Value = "06:52:41"
$ws.Cells.Item(23,2).Value = "06:56"
$ws.Cells.Item(23,3).Value = "14_ABASTO"
$ws.Cells.Item(23,4).Value = 4
$ws.Cells.Item(23,5).Value = "LP1912"
$ws.Cells.Item(24,1).Value = "06:03:38"
$ws.Cells.Item(24,2).Value = "06:59"
$ws.Cells.Item(24,3).Value = "14_ABASTO"
$ws.Cells.Item(24,4).Value = 56
$ws.Cells.Item(24,5).Value = "LP1912"
$ws.Cells.Item(25,1).Value = "06:52:41"
$ws.Cells.Item(25,2).Value = "07:01"
$ws.Cells.Item(25,3).Value = "16_SANTA ANA"
$ws.Cells.Item(25,4).Value = 9
$ws.Cells.Item(25,5).Value = "LP1912"
$ws.Cells.Item(26,1).Value = "06:52:41"
$ws.Cells.Item(26,2).Value = "07:04"
$ws.Cells.Item(26,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(26,4).Value = 12
$ws.Cells.Item(26,5).Value = "LP1912"
$ws.Cells.Item(27,1).Value = "06:52:41"
$ws.Cells.Item(27,2).Value = "07:05"
$ws.Cells.Item(27,3).Value = "15_ABASTO"
$ws.Cells.Item(27,4).Value = 13
$ws.Cells.Item(27,5).Value = "LP1912"
$ws.Cells.Item(28,1).Value = "06:52:41"
$ws.Cells.Item(28,2).Value = "07:07"
$ws.Cells.Item(28,3).Value = "225_GOMEZ"
$ws.Cells.Item(28,4).Value = 15
$ws.Cells.Item(28,5).Value = "LP1912"
$ws.Cells.Item(29,1).Value = "06:52:41"
$ws.Cells.Item(29,2).Value = "07:11"
$ws.Cells.Item(29,3).Value = "215A_EL PATO"
$ws.Cells.Item(29,4).Value = 19
$ws.Cells.Item(29,5).Value = "LP1912"
$ws.Cells.Item(30,1).Value = "06:52:41"
$ws.Cells.Item(30,2).Value = "07:15"
$ws.Cells.Item(30,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(30,4).Value = 23
$ws.Cells.Item(30,5).Value = "LP1912"
$ws.Cells.Item(31,1).Value = "06:52:41"
$ws.Cells.Item(31,2).Value = "07:16"
$ws.Cells.Item(31,3).Value = "16_SANTA ANA"
$ws.Cells.Item(31,4).Value = 24
$ws.Cells.Item(31,5).Value = "LP1912"
$ws.Cells.Item(32,1).Value = "06:52:41"
$ws.Cells.Item(32,2).Value = "07:21"
$ws.Cells.Item(32,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(32,4).Value = 29
$ws.Cells.Item(32,5).Value = "LP1912"
$ws.Cells.Item(33,1).Value = "07:23:38"
$ws.Cells.Item(33,2).Value = "07:23"
$ws.Cells.Item(33,3).Value = "16_SANTA ANA"
$ws.Cells.Item(33,4).Value = 0
$ws.Cells.Item(33,5).Value = "LP1912"
$ws.Cells.Item(34,1).Value = "06:52:41"
$ws.Cells.Item(34,2).Value = "07:23"
$ws.Cells.Item(34,3).Value = "10_OLMOS"
$ws.Cells.Item(34,4).Value = 31
$ws.Cells.Item(34,5).Value = "LP1912"
$ws.Cells.Item(35,1).Value = "07:23:38"
$ws.Cells.Item(35,2).Value = "07:31"
$ws.Cells.Item(35,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(35,4).Value = 8
$ws.Cells.Item(35,5).Value = "LP1912"
$ws.Cells.Item(36,1).Value = "07:23:38"
$ws.Cells.Item(36,2).Value = "07:31"
$ws.Cells.Item(36,3).Value = "16_SANTA ANA"
$ws.Cells.Item(36,4).Value = 8
$ws.Cells.Item(36,5).Value = "LP1912"
$ws.Cells.Item(37,1).Value = "07:23:38"
$ws.Cells.Item(37,2).Value = "07:32"
$ws.Cells.Item(37,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(37,4).Value = 9
$ws.Cells.Item(37,5).Value = "LP1912"
$ws.Cells.Item(38,1).Value = "07:23:38"
$ws.Cells.Item(38,2).Value = "07:36"
$ws.Cells.Item(38,3).Value = "27_EL RETIRO"
$ws.Cells.Item(38,4).Value = 13
$ws.Cells.Item(38,5).Value = "LP1912"
$ws.Cells.Item(39,1).Value = "07:23:38"
$ws.Cells.Item(39,2).Value = "07:36"
$ws.Cells.Item(39,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(39,4).Value = 13
$ws.Cells.Item(39,5).Value = "LP1912"
$ws.Cells.Item(40,1).Value = "07:23:38"
$ws.Cells.Item(40,2).Value = "07:39"
$ws.Cells.Item(40,3).Value = "10_OLMOS"
$ws.Cells.Item(40,4).Value = 16
$ws.Cells.Item(40,5).Value = "LP1912"
$ws.Cells.Item(41,1).Value = "07:23:38"
$ws.Cells.Item(41,2).Value = "07:47"
$ws.Cells.Item(41,3).Value = "14_ABASTO"
$ws.Cells.Item(41,4).Value = 24
$ws.Cells.Item(41,5).Value = "LP1912"
$ws.Cells.Item(42,1).Value = "07:23:38"
$ws.Cells.Item(42,2).Value = "07:51"
$ws.Cells.Item(42,3).Value = "215D_EL PATO"
$ws.Cells.Item(42,4).Value = 28
$ws.Cells.Item(42,5).Value = "LP1912"
$ws.Cells.Item(43,1).Value = "07:23:38"
$ws.Cells.Item(43,2).Value = "07:55"
$ws.Cells.Item(43,3).Value = "10_OLMOS"
$ws.Cells.Item(43,4).Value = 32
$ws.Cells.Item(43,5).Value = "LP1912"
$ws.Cells.Item(44,1).Value = "07:59:51"
$ws.Cells.Item(44,2).Value = "07:59"
$ws.Cells.Item(44,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = "LP1912"
$ws.Cells.Item(45,1).Value = "07:59:51"
$ws.Cells.Item(45,2).Value = "07:59"
$ws.Cells.Item(45,3).Value = "16_SANTA ANA"
$ws.Cells.Item(45,4).Value = 0
$ws.Cells.Item(45,5).Value = "LP1912"
$ws.Cells.Item(46,1).Value = "07:23:38"
$ws.Cells.Item(46,2).Value = "08:03"
$ws.Cells.Item(46,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(46,4).Value = 40
$ws.Cells.Item(46,5).Value = "LP1912"
$ws.Cells.Item(47,1).Value = "07:59:51"
$ws.Cells.Item(47,2).Value = "08:04"
$ws.Cells.Item(47,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(47,4).Value = 5
$ws.Cells.Item(47,5).Value = "LP1912"
$ws.Cells.Item(48,1).Value = "06:52:41"
$ws.Cells.Item(48,2).Value = "08:06"
$ws.Cells.Item(48,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(48,4).Value = 74
$ws.Cells.Item(48,5).Value = "LP1912"
$ws.Cells.Item(49,1).Value = "07:59:51"
$ws.Cells.Item(49,2).Value = "08:12"
$ws.Cells.Item(49,3).Value = "15_ABASTO"
$ws.Cells.Item(49,4).Value = 13
$ws.Cells.Item(49,5).Value = "LP1912"
$ws.Cells.Item(50,1).Value = "07:59:51"
$ws.Cells.Item(50,2).Value = "08:13"
$ws.Cells.Item(50,3).Value = "10_OLMOS"
$ws.Cells.Item(50,4).Value = 14
$ws.Cells.Item(50,5).Value = "LP1912"
$ws.Cells.Item(51,1).Value = "07:59:51"
$ws.Cells.Item(51,2).Value = "08:21"
$ws.Cells.Item(51,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(51,4).Value = 22
$ws.Cells.Item(51,5).Value = "LP1912"
$ws.Cells.Item(52,1).Value = "07:59:51"
$ws.Cells.Item(52,2).Value = "08:22"
$ws.Cells.Item(52,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(52,4).Value = 23
$ws.Cells.Item(52,5).Value = "LP1912"
$ws.Cells.Item(53,1).Value = "07:59:51"
$ws.Cells.Item(53,2).Value = "08:23"
$ws.Cells.Item(53,3).Value = "215B_EL PATO"
$ws.Cells.Item(53,4).Value = 24
$ws.Cells.Item(53,5).Value = "LP1912"
$ws.Cells.Item(54,1).Value = "07:59:51"
$ws.Cells.Item(54,2).Value = "08:27"
$ws.Cells.Item(54,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(54,4).Value = 28
$ws.Cells.Item(54,5).Value = "LP1912"
$ws.Cells.Item(55,1).Value = "08:31:16"
$ws.Cells.Item(55,2).Value = "08:32"
$ws.Cells.Item(55,3).Value = "10_OLMOS"
$ws.Cells.Item(55,4).Value = 1
$ws.Cells.Item(55,5).Value = "LP1912"
$ws.Cells.Item(56,1).Value = "07:59:51"
$ws.Cells.Item(56,2).Value = "08:34"
$ws.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56,4).Value = 35
$ws.Cells.Item(56,5).Value = "LP1912"
$ws.Cells.Item(57,1).Value = "08:31:16"
$ws.Cells.Item(57,2).Value = "08:35"
$ws.Cells.Item(57,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(57,4).Value = 4
$ws.Cells.Item(57,5).Value = "LP1912"
$ws.Cells.Item(58,1).Value = "08:31:16"
$ws.Cells.Item(58,2).Value = "08:42"
$ws.Cells.Item(58,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(58,4).Value = 11
$ws.Cells.Item(58,5).Value = "LP1912"
$ws.Cells.Item(59,1).Value = "07:59:51"
$ws.Cells.Item(59,2).Value = "08:43"
$ws.Cells.Item(59,3).Value = "14_ABASTO"
$ws.Cells.Item(59,4).Value = 44
$ws.Cells.Item(59,5).Value = "LP1912"
$ws.Cells.Item(60,1).Value = "08:31:16"
$ws.Cells.Item(60,2).Value = "08:44"
$ws.Cells.Item(60,3).Value = "14_ABASTO"
$ws.Cells.Item(60,4).Value = 13
$ws.Cells.Item(60,5).Value = "LP1912"
$ws.Cells.Item(61,1).Value = "07:23:38"
$ws.Cells.Item(61,2).Value = "08:44"
$ws.Cells.Item(61,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(61,4).Value = 81
$ws.Cells.Item(61,5).Value = "LP1912"
$ws.Cells.Item(62,1).Value = "08:31:16"
$ws.Cells.Item(62,2).Value = "08:53"
$ws.Cells.Item(62,3).Value = "10_OLMOS"
$ws.Cells.Item(62,4).Value = 22
$ws.Cells.Item(62,5).Value = "LP1912"
$ws.Cells.Item(63,1).Value = "08:54:41"
$ws.Cells.Item(63,2).Value = "08:54"
$ws.Cells.Item(63,3).Value = "17_ROMERO"
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = "LP1912"
$ws.Cells.Item(64,1).Value = "08:54:41"
$ws.Cells.Item(64,2).Value = "08:57"
$ws.Cells.Item(64,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(64,4).Value = 3
$ws.Cells.Item(64,5).Value = "LP1912"
$ws.Cells.Item(65,1).Value = "08:31:16"
$ws.Cells.Item(65,2).Value = "09:01"
$ws.Cells.Item(65,3).Value = "215A_EL PATO"
$ws.Cells.Item(65,4).Value = 30
$ws.Cells.Item(65,5).Value = "LP1912"
$ws.Cells.Item(66,1).Value = "08:54:41"
$ws.Cells.Item(66,2).Value = "09:03"
$ws.Cells.Item(66,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(66,4).Value = 9
$ws.Cells.Item(66,5).Value = "LP1912"
$ws.Cells.Item(67,1).Value = "08:31:16"
$ws.Cells.Item(67,2).Value = "09:04"
$ws.Cells.Item(67,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(67,4).Value = 33
$ws.Cells.Item(67,5).Value = "LP1912"
$ws.Cells.Item(68,1).Value = "08:31:16"
$ws.Cells.Item(68,2).Value = "09:05"
$ws.Cells.Item(68,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(68,4).Value = 34
$ws.Cells.Item(68,5).Value = "LP1912"
$ws.Cells.Item(69,1).Value = "08:54:41"
$ws.Cells.Item(69,2).Value = "09:06"
$ws.Cells.Item(69,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(69,4).Value = 12
$ws.Cells.Item(69,5).Value = "LP1912"
$ws.Cells.Item(70,1).Value = "08:54:41"
$ws.Cells.Item(70,2).Value = "09:10"
$ws.Cells.Item(70,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(70,4).Value = 16
$ws.Cells.Item(70,5).Value = "LP1912"
$ws.Cells.Item(71,1).Value = "08:31:16"
$ws.Cells.Item(71,2).Value = "09:11"
$ws.Cells.Item(71,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(71,4).Value = 40
$ws.Cells.Item(71,5).Value = "LP1912"
$ws.Cells.Item(72,1).Value = "08:54:41"
$ws.Cells.Item(72,2).Value = "09:13"
$ws.Cells.Item(72,3).Value = "10_OLMOS"
$ws.Cells.Item(72,4).Value = 19
$ws.Cells.Item(72,5).Value = "LP1912"
$ws.Cells.Item(73,1).Value = "08:54:41"
$ws.Cells.Item(73,2).Value = "09:16"
$ws.Cells.Item(73,3).Value = "27_EL RETIRO"
$ws.Cells.Item(73,4).Value = 22
$ws.Cells.Item(73,5).Value = "LP1912"
$ws.Cells.Item(74,1).Value = "08:31:16"
$ws.Cells.Item(74,2).Value = "09:17"
$ws.Cells.Item(74,3).Value = "27_EL RETIRO"
$ws.Cells.Item(74,4).Value = 46
$ws.Cells.Item(74,5).Value = "LP1912"
$ws.Cells.Item(75,1).Value = "07:59:51"
$ws.Cells.Item(75,2).Value = "09:20"
$ws.Cells.Item(75,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(75,4).Value = 81
$ws.Cells.Item(75,5).Value = "LP1912"
$ws.Cells.Item(76,1).Value = "08:54:41"
$ws.Cells.Item(76,2).Value = "09:21"
$ws.Cells.Item(76,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(76,4).Value = 27
$ws.Cells.Item(76,5).Value = "LP1912"
$ws.Cells.Item(77,1).Value = "07:59:51"
$ws.Cells.Item(77,2).Value = "09:22"
$ws.Cells.Item(77,3).Value = "17_ROMERO"
$ws.Cells.Item(77,4).Value = 83
$ws.Cells.Item(77,5).Value = "LP1912"
$ws.Cells.Item(78,1).Value = "08:54:41"
$ws.Cells.Item(78,2).Value = "09:22"
$ws.Cells.Item(78,3).Value = "16_SANTA ANA"
$ws.Cells.Item(78,4).Value = 28
$ws.Cells.Item(78,5).Value = "LP1912"
$ws.Cells.Item(79,1).Value = "08:31:16"
$ws.Cells.Item(79,2).Value = "09:23"
$ws.Cells.Item(79,3).Value = "16_SANTA ANA"
$ws.Cells.Item(79,4).Value = 52
$ws.Cells.Item(79,5).Value = "LP1912"
$ws.Cells.Item(80,1).Value = "08:54:41"
$ws.Cells.Item(80,2).Value = "09:23"
$ws.Cells.Item(80,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(80,4).Value = 29
$ws.Cells.Item(80,5).Value = "LP1912"
$ws.Cells.Item(81,1).Value = "08:31:16"
$ws.Cells.Item(81,2).Value = "09:24"
$ws.Cells.Item(81,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(81,4).Value = 53
$ws.Cells.Item(81,5).Value = "LP1912"
$ws.Cells.Item(82,1).Value = "08:54:41"
$ws.Cells.Item(82,2).Value = "09:32"
$ws.Cells.Item(82,3).Value = "15_ABASTO"
$ws.Cells.Item(82,4).Value = 38
$ws.Cells.Item(82,5).Value = "LP1912"
$ws.Cells.Item(83,1).Value = "09:32:47"
$ws.Cells.Item(83,2).Value = "09:33"
$ws.Cells.Item(83,3).Value = "10_OLMOS"
$ws.Cells.Item(83,4).Value = 1
$ws.Cells.Item(83,5).Value = "LP1912"
$ws.Cells.Item(84,1).Value = "09:32:47"
$ws.Cells.Item(84,2).Value = "09:34"
$ws.Cells.Item(84,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(84,4).Value = 2
$ws.Cells.Item(84,5).Value = "LP1912"
$ws.Cells.Item(85,1).Value = "08:54:41"
$ws.Cells.Item(85,2).Value = "09:34"
$ws.Cells.Item(85,3).Value = "16_SANTA ANA"
$ws.Cells.Item(85,4).Value = 40
$ws.Cells.Item(85,5).Value = "LP1912"
$ws.Cells.Item(86,1).Value = "08:31:16"
$ws.Cells.Item(86,2).Value = "09:35"
$ws.Cells.Item(86,3).Value = "16_SANTA ANA"
$ws.Cells.Item(86,4).Value = 64
$ws.Cells.Item(86,5).Value = "LP1912"
$ws.Cells.Item(87,1).Value = "09:32:47"
$ws.Cells.Item(87,2).Value = "09:41"
$ws.Cells.Item(87,3).Value = "215C_EL PATO"
$ws.Cells.Item(87,4).Value = 9
$ws.Cells.Item(87,5).Value = "LP1912"
$ws.Cells.Item(88,1).Value = "08:31:16"
$ws.Cells.Item(88,2).Value = "09:42"
$ws.Cells.Item(88,3).Value = "215C_EL PATO"
$ws.Cells.Item(88,4).Value = 71
$ws.Cells.Item(88,5).Value = "LP1912"
$ws.Cells.Item(89,1).Value = "09:32:47"
$ws.Cells.Item(89,2).Value = "09:42"
$ws.Cells.Item(89,3).Value = "16_SANTA ANA"
$ws.Cells.Item(89,4).Value = 10
$ws.Cells.Item(89,5).Value = "LP1912"
$ws.Cells.Item(90,1).Value = "09:32:47"
$ws.Cells.Item(90,2).Value = "09:43"
$ws.Cells.Item(90,3).Value = "14_ABASTO"
$ws.Cells.Item(90,4).Value = 11
$ws.Cells.Item(90,5).Value = "LP1912"
$ws.Cells.Item(91,1).Value = "08:31:16"
$ws.Cells.Item(91,2).Value = "09:44"
$ws.Cells.Item(91,3).Value = "14_ABASTO"
$ws.Cells.Item(91,4).Value = 73
$ws.Cells.Item(91,5).Value = "LP1912"
$ws.Cells.Item(92,1).Value = "09:32:47"
$ws.Cells.Item(92,2).Value = "09:46"
$ws.Cells.Item(92,3).Value = "16_SANTA ANA"
$ws.Cells.Item(92,4).Value = 14
$ws.Cells.Item(92,5).Value = "LP1912"
$ws.Cells.Item(93,1).Value = "08:54:41"
$ws.Cells.Item(93,2).Value = "09:52"
$ws.Cells.Item(93,3).Value = "15_ABASTO"
$ws.Cells.Item(93,4).Value = 58
$ws.Cells.Item(93,5).Value = "LP1912"
$ws.Cells.Item(94,1).Value = "09:32:47"
$ws.Cells.Item(94,2).Value = "09:53"
$ws.Cells.Item(94,3).Value = "10_OLMOS"
$ws.Cells.Item(94,4).Value = 21
$ws.Cells.Item(94,5).Value = "LP1912"
$ws.Cells.Item(95,1).Value = "09:32:47"
$ws.Cells.Item(95,2).Value = "09:54"
$ws.Cells.Item(95,3).Value = "15_ABASTO"
$ws.Cells.Item(95,4).Value = 22
$ws.Cells.Item(95,5).Value = "LP1912"
$ws.Cells.Item(96,1).Value = "09:32:47"
$ws.Cells.Item(96,2).Value = "10:03"
$ws.Cells.Item(96,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(96,4).Value = 31
$ws.Cells.Item(96,5).Value = "LP1912"
$ws.Cells.Item(97,1).Value = "09:32:47"
$ws.Cells.Item(97,2).Value = "10:04"
$ws.Cells.Item(97,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(97,4).Value = 32
$ws.Cells.Item(97,5).Value = "LP1912"
$ws.Cells.Item(98,1).Value = "09:32:47"
$ws.Cells.Item(98,2).Value = "10:12"
$ws.Cells.Item(98,3).Value = "15_ABASTO"
$ws.Cells.Item(98,4).Value = 40
$ws.Cells.Item(98,5).Value = "LP1912"
$ws.Cells.Item(99,1).Value = "09:32:47"
$ws.Cells.Item(99,2).Value = "10:13"
$ws.Cells.Item(99,3).Value = "10_OLMOS"
$ws.Cells.Item(99,4).Value = 41
$ws.Cells.Item(99,5).Value = "LP1912"
$ws.Cells.Item(100,1).Value = "09:32:47"
$ws.Cells.Item(100,2).Value = "10:21"
$ws.Cells.Item(100,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(100,4).Value = 49
$ws.Cells.Item(100,5).Value = "LP1912"
$ws.Cells.Item(101,1).Value = "09:32:47"
$ws.Cells.Item(101,2).Value = "10:22"
$ws.Cells.Item(101,3).Value = "17_ROMERO"
$ws.Cells.Item(101,4).Value = 50
$ws.Cells.Item(101,5).Value = "LP1912"
$ws.Cells.Item(102,1).Value = "09:32:47"
$ws.Cells.Item(102,2).Value = "10:23"
$ws.Cells.Item(102,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(102,4).Value = 51
$ws.Cells.Item(102,5).Value = "LP1912"
$ws.Cells.Item(103,1).Value = "09:32:47"
$ws.Cells.Item(103,2).Value = "10:26"
$ws.Cells.Item(103,3).Value = "215A_EL PATO"
$ws.Cells.Item(103,4).Value = 54
$ws.Cells.Item(103,5).Value = "LP1912"
$ws.Cells.Item(104,1).Value = "09:32:47"
$ws.Cells.Item(104,2).Value = "10:41"
$ws.Cells.Item(104,3).Value = "17_ROMERO"
$ws.Cells.Item(104,4).Value = 69
$ws.Cells.Item(104,5).Value = "LP1912"
$ws.Cells.Item(105,1).Value = "09:32:47"
$ws.Cells.Item(105,2).Value = "10:43"
$ws.Cells.Item(105,3).Value = "14_ABASTO"
$ws.Cells.Item(105,4).Value = 71
$ws.Cells.Item(105,5).Value = "LP1912"
$ws.Cells.Item(106,1).Value = "09:32:47"
$ws.Cells.Item(106,2).Value = "10:53"
$ws.Cells.Item(106,3).Value = "27_EL RETIRO"
$ws.Cells.Item(106,4).Value = 81
$ws.Cells.Item(106,5).Value = "LP1912"
$ws.Cells.Item(107,1).Value = "09:32:47"
$ws.Cells.Item(107,2).Value = "11:01"
$ws.Cells.Item(107,3).Value = "215C_EL PATO"
$ws.Cells.Item(107,4).Value = 89
$ws.Cells.Item(107,5).Value = "LP1912"
$ws.Cells.Item(108,1).Value = "09:32:47"
$ws.Cells.Item(108,2).Value = "11:06"
$ws.Cells.Item(108,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(108,4).Value = 94
$ws.Cells.Item(108,5).Value = "LP1912"
$ws.Cells.Item(109,1).Value = "09:32:47"
$ws.Cells.Item(109,2).Value = "11:19"
$ws.Cells.Item(109,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(109,4).Value = 107
$ws.Cells.Item(109,5).Value = "LP1912"
$ws.Cells.Item(110,1).Value = "09:32:47"
$ws.Cells.Item(110,2).Value = "11:21"
$ws.Cells.Item(110,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(110,4).Value = 109
$ws.Cells.Item(110,5).Value = "LP1912"
$ws.Cells.Item(111,1).Value = "09:32:47"
$ws.Cells.Item(111,2).Value = "11:26"
$ws.Cells.Item(111,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(111,4).Value = 114
$ws.Cells.Item(111,5).Value = "LP1912"
$ws.Cells.Item(112,1).Value = "09:32:47"
$ws.Cells.Item(112,2).Value = "11:27"
$ws.Cells.Item(112,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(112,4).Value = 115
$ws.Cells.Item(112,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 09:32:47"
$ws.Cells.Item(3,1).Value = "Total filas: 12"
$ws.Cells.Item(6,1).Value = "04:03:21"
$ws.Cells.Item(6,2).Value = "04:46"
$ws.Cells.Item(6,3).Value = "215A_EL PATO"
$ws.Cells.Item(6,4).Value = 43
$ws.Cells.Item(6,5).Value = "LP1912"
$ws.Cells.Item(7,1).Value = "05:20:44"
$ws.Cells.Item(7,2).Value = "05:34"
$ws.Cells.Item(7,3).Value = "215B_EL PATO"
$ws.Cells.Item(7,4).Value = 14
$ws.Cells.Item(7,5).Value = "LP1912"
$ws.Cells.Item(8,1).Value = "06:03:38"
$ws.Cells.Item(8,2).Value = "06:11"
$ws.Cells.Item(8,3).Value = "215A_EL PATO"
$ws.Cells.Item(8,4).Value = 8
$ws.Cells.Item(8,5).Value = "LP1912"
$ws.Cells.Item(9,1).Value = "06:03:38"
$ws.Cells.Item(9,2).Value = "06:46"
$ws.Cells.Item(9,3).Value = "215C_EL PATO"
$ws.Cells.Item(9,4).Value = 43
$ws.Cells.Item(9,5).Value = "LP1912"
$ws.Cells.Item(10,1).Value = "06:52:41"
$ws.Cells.Item(10,2).Value = "07:11"
$ws.Cells.Item(10,3).Value = "215A_EL PATO"
$ws.Cells.Item(10,4).Value = 19
$ws.Cells.Item(10,5).Value = "LP1912"
$ws.Cells.Item(11,1).Value = "07:23:38"
$ws.Cells.Item(11,2).Value = "07:51"
$ws.Cells.Item(11,3).Value = "215D_EL PATO"
$ws.Cells.Item(11,4).Value = 28
$ws.Cells.Item(11,5).Value = "LP1912"
$ws.Cells.Item(12,1).Value = "07:59:51"
$ws.Cells.Item(12,2).Value = "08:23"
$ws.Cells.Item(12,3).Value = "215B_EL PATO"
$ws.Cells.Item(12,4).Value = 24
$ws.Cells.Item(12,5).Value = "LP1912"
$ws.Cells.Item(13,1).Value = "08:31:16"
$ws.Cells.Item(13,2).Value = "09:01"
$ws.Cells.Item(13,3).Value = "215A_EL PATO"
$ws.Cells.Item(13,4).Value = 30
$ws.Cells.Item(13,5).Value = "LP1912"
$ws.Cells.Item(14,1).Value = "09:32:47"
$ws.Cells.Item(14,2).Value = "09:41"
$ws.Cells.Item(14,3).Value = "215C_EL PATO"
$ws.Cells.Item(14,4).Value = 9
$ws.Cells.Item(14,5).Value = "LP1912"
$ws.Cells.Item(15,1).Value = "08:31:16"
$ws.Cells.Item(15,2).Value = "09:42"
$ws.Cells.Item(15,3).Value = "215C_EL PATO"
$ws.Cells.Item(15,4).Value = 71
$ws.Cells.Item(15,5).Value = "LP1912"
$ws.Cells.Item(16,1).Value = "09:32:47"
$ws.Cells.Item(16,2).Value = "10:26"
$ws.Cells.Item(16,3).Value = "215A_EL PATO"
$ws.Cells.Item(16,4).Value = 54
$ws.Cells.Item(16,5).Value = "LP1912"
$ws.Cells.Item(17,1).Value = "09:32:47"
$ws.Cells.Item(17,2).Value = "11:01"
$ws.Cells.Item(17,3).Value = "215C_EL PATO"
$ws.Cells.Item(17,4).Value = 89
$ws.Cells.Item(17,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 09:32:47"
$ws.Cells.Item(3,1).Value = "Total filas: 19"
$ws.Cells.Item(6,1).Value = "05:20:44"
$ws.Cells.Item(6,2).Value = "05:43"
$ws.Cells.Item(6,3).Value = "215A_LA PLATA"
$ws.Cells.Item(6,4).Value = 23
$ws.Cells.Item(6,5).Value = "L6173"
$ws.Cells.Item(7,1).Value = "04:03:21"
$ws.Cells.Item(7,2).Value = "05:44"
$ws.Cells.Item(7,3).Value = "215A_LA PLATA"
$ws.Cells.Item(7,4).Value = 101
$ws.Cells.Item(7,5).Value = "L6173"
$ws.Cells.Item(8,1).Value = "06:03:38"
$ws.Cells.Item(8,2).Value = "06:08"
$ws.Cells.Item(8,3).Value = "215A_LA PLATA"
$ws.Cells.Item(8,4).Value = 5
$ws.Cells.Item(8,5).Value = "L6173"
$ws.Cells.Item(9,1).Value = "06:03:38"
$ws.Cells.Item(9,2).Value = "06:32"
$ws.Cells.Item(9,3).Value = "215C_LA PLATA"
$ws.Cells.Item(9,4).Value = 29
$ws.Cells.Item(9,5).Value = "L6203"
$ws.Cells.Item(10,1).Value = "06:52:41"
$ws.Cells.Item(10,2).Value = "06:59"
$ws.Cells.Item(10,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(10,4).Value = 7
$ws.Cells.Item(10,5).Value = "L6173"
$ws.Cells.Item(11,1).Value = "06:03:38"
$ws.Cells.Item(11,2).Value = "07:00"
$ws.Cells.Item(11,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(11,4).Value = 57
$ws.Cells.Item(11,5).Value = "L6173"
$ws.Cells.Item(12,1).Value = "07:23:38"
$ws.Cells.Item(12,2).Value = "07:35"
$ws.Cells.Item(12,3).Value = "215A_LA PLATA"
$ws.Cells.Item(12,4).Value = 12
$ws.Cells.Item(12,5).Value = "L6173"
$ws.Cells.Item(13,1).Value = "06:52:41"
$ws.Cells.Item(13,2).Value = "08:06"
$ws.Cells.Item(13,3).Value = "215C_LA PLATA"
$ws.Cells.Item(13,4).Value = 74
$ws.Cells.Item(13,5).Value = "L6203"
$ws.Cells.Item(14,1).Value = "07:23:38"
$ws.Cells.Item(14,2).Value = "08:09"
$ws.Cells.Item(14,3).Value = "215C_LA PLATA"
$ws.Cells.Item(14,4).Value = 46
$ws.Cells.Item(14,5).Value = "L6203"
$ws.Cells.Item(15,1).Value = "07:59:51"
$ws.Cells.Item(15,2).Value = "08:11"
$ws.Cells.Item(15,3).Value = "215C_LA PLATA"
$ws.Cells.Item(15,4).Value = 12
$ws.Cells.Item(15,5).Value = "L6203"
$ws.Cells.Item(16,1).Value = "06:52:41"
$ws.Cells.Item(16,2).Value = "08:31"
$ws.Cells.Item(16,3).Value = "215A_LA PLATA"
$ws.Cells.Item(16,4).Value = 99
$ws.Cells.Item(16,5).Value = "L6173"
$ws.Cells.Item(17,1).Value = "07:59:51"
$ws.Cells.Item(17,2).Value = "08:35"
$ws.Cells.Item(17,3).Value = "215A_LA PLATA"
$ws.Cells.Item(17,4).Value = 36
$ws.Cells.Item(17,5).Value = "L6173"
$ws.Cells.Item(18,1).Value = "08:31:16"
$ws.Cells.Item(18,2).Value = "08:37"
$ws.Cells.Item(18,3).Value = "215A_LA PLATA"
$ws.Cells.Item(18,4).Value = 6
$ws.Cells.Item(18,5).Value = "L6173"
$ws.Cells.Item(19,1).Value = "07:59:51"
$ws.Cells.Item(19,2).Value = "09:08"
$ws.Cells.Item(19,3).Value = "215D_LA PLATA"
$ws.Cells.Item(19,4).Value = 69
$ws.Cells.Item(19,5).Value = "L6203"
$ws.Cells.Item(20,1).Value = "08:54:41"
$ws.Cells.Item(20,2).Value = "09:09"
$ws.Cells.Item(20,3).Value = "215D_LA PLATA"
$ws.Cells.Item(20,4).Value = 15
$ws.Cells.Item(20,5).Value = "L6203"
$ws.Cells.Item(21,1).Value = "09:32:47"
$ws.Cells.Item(21,2).Value = "10:02"
$ws.Cells.Item(21,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(21,4).Value = 30
$ws.Cells.Item(21,5).Value = "L6173"
$ws.Cells.Item(22,1).Value = "08:31:16"
$ws.Cells.Item(22,2).Value = "10:03"
$ws.Cells.Item(22,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(22,4).Value = 92
$ws.Cells.Item(22,5).Value = "L6173"
$ws.Cells.Item(23,1).Value = "09:32:47"
$ws.Cells.Item(23,2).Value = "10:54"
$ws.Cells.Item(23,3).Value = "215A_LA PLATA"
$ws.Cells.Item(23,4).Value = 82
$ws.Cells.Item(23,5).Value = "L6173"
$ws.Cells.Item(24,1).Value = "09:32:47"
$ws.Cells.Item(24,2).Value = "11:13"
$ws.Cells.Item(24,3).Value = "215C_LA PLATA"
$ws.Cells.Item(24,4).Value = 101
$ws.Cells.Item(24,5).Value = "L6203"
